$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Add new "Email (T/F)" column (G) - header + data value, matching
#    the formatting already used by the neighbouring "Website (T/F)"
#    column (F).
# ------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Email (T/F)"

$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "F"

# ------------------------------------------------------------------
# 2. Extend the generated "mod_info_server(...)" formula in B6 so it
#    also emits the new email = ... argument.
# ------------------------------------------------------------------
$ws.Range("B6").Formula = "=""mod_info_server('""&A2&""', selector = selection, data = ""&C2&"", rownametitle = c('""&B2&""'), phone = ""&E2&"", website = ""&F2&"", email = ""&G2&"")"""

# ------------------------------------------------------------------
# 3. Remove the old ICA reference/hyperlink scratch row (row 9) that
#    held the resource's phone/website/email/description -- this is
#    no longer needed now that the row-2 table drives the formulas.
#    Deleting the whole row shifts rows 10-14 up to 9-13.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Rows(9).Delete()

# ------------------------------------------------------------------
# 4. Clear the leftover "ENTER INTO UI" helper labels that used to sit
#    in column A of the (now renumbered) rows 11-13.
# ------------------------------------------------------------------
$ws.Range("A11").Clear()
$ws.Range("A12").Clear()
$ws.Range("A13").Clear()

# ------------------------------------------------------------------
# 5. Update the view: scroll so column B is left-most and select G6,
#    mirroring where the author was working when they saved.
# ------------------------------------------------------------------
$ws.Range("G6").Select()

Write-Host "edit complete"
